$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), copying the format/style
# from the existing header cell H1 so they share the same style index.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I ("I0") and J ("IF"), rows 2-9.
$data = @{
    2 = @(1, 5)
    3 = @(1, 6)
    4 = @(1, 5)
    5 = @(1, 4)
    6 = @(6, 7)
    7 = @(1, 2)
    8 = @(1, 2)
    9 = @(7, 7)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
